$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Cxcl13"
$ws.Range("C2").Value = "Cxcr5"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.009847666666666
$ws.Range("H2").Value = 12.029543
$ws.Range("I2").Value = 0.9697248931871538
$ws.Range("J2").Value = 0.9697248931871538
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.6511963333333334
$ws.Range("N2").Value = 1.953589
$ws.Range("O2").Value = 0.3942798821674536
$ws.Range("P2").Value = 0.3942798821674536
$ws.Range("Q2").Value = 2.611198097758556
$ws.Range("R2").Value = 23.500782879827
$ws.Range("S2").Value = 0.3823430166206775
$ws.Range("T2").Value = 0.3823430166206775

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Cxcl13"
$ws.Range("C3").Value = "Cxcr5"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.009847666666666
$ws.Range("H3").Value = 12.029543
$ws.Range("I3").Value = 0.9697248931871538
$ws.Range("J3").Value = 0.9697248931871538
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.6731493333333334
$ws.Range("N3").Value = 2.019448
$ws.Range("O3").Value = 0.4075717663660575
$ws.Range("P3").Value = 0.4075717663660575
$ws.Range("Q3").Value = 2.699226283584889
$ws.Range("R3").Value = 24.293036552264
$ws.Range("S3").Value = 0.3952324876054247
$ws.Range("T3").Value = 0.3952324876054247

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cxcl13"
$ws.Range("C4").Value = "Cxcr5"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.009847666666666
$ws.Range("H4").Value = 12.029543
$ws.Range("I4").Value = 0.9697248931871538
$ws.Range("J4").Value = 0.9697248931871538
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3272636666666667
$ws.Range("N4").Value = 0.9817910000000001
$ws.Range("O4").Value = 0.1981483514664888
$ws.Range("P4").Value = 0.1981483514664889
$ws.Range("Q4").Value = 1.312277450168111
$ws.Range("R4").Value = 11.810497051513
$ws.Range("S4").Value = 0.1921493889610515
$ws.Range("T4").Value = 0.1921493889610515

# Row 5
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Cxcl13"
$ws.Range("C5").Value = "Cxcr5"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1251886666666667
$ws.Range("H5").Value = 0.375566
$ws.Range("I5").Value = 0.03027510681284623
$ws.Range("J5").Value = 0.03027510681284622
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.6511963333333334
$ws.Range("N5").Value = 1.953589
$ws.Range("O5").Value = 0.3942798821674536
$ws.Range("P5").Value = 0.3942798821674536
$ws.Range("Q5").Value = 0.08152240070822223
$ws.Range("R5").Value = 0.7337016063740001
$ws.Range("S5").Value = 0.01193686554677608
$ws.Range("T5").Value = 0.01193686554677608

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Cxcl13"
$ws.Range("C6").Value = "Cxcr5"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1251886666666667
$ws.Range("H6").Value = 0.375566
$ws.Range("I6").Value = 0.03027510681284623
$ws.Range("J6").Value = 0.03027510681284622
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.6731493333333334
$ws.Range("N6").Value = 2.019448
$ws.Range("O6").Value = 0.4075717663660575
$ws.Range("P6").Value = 0.4075717663660575
$ws.Range("Q6").Value = 0.08427066750755556
$ws.Range("R6").Value = 0.7584360075680001
$ws.Range("S6").Value = 0.0123392787606328
$ws.Range("T6").Value = 0.0123392787606328

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Cxcl13"
$ws.Range("C7").Value = "Cxcr5"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.1251886666666667
$ws.Range("H7").Value = 0.375566
$ws.Range("I7").Value = 0.03027510681284623
$ws.Range("J7").Value = 0.03027510681284622
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.3272636666666667
$ws.Range("N7").Value = 0.9817910000000001
$ws.Range("O7").Value = 0.1981483514664888
$ws.Range("P7").Value = 0.1981483514664889
$ws.Range("Q7").Value = 0.04096970207844445
$ws.Range("R7").Value = 0.3687273187060001
$ws.Range("S7").Value = 0.005998962505437345
$ws.Range("T7").Value = 0.005998962505437346

